$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 171.0
$ws.Range("B3").Value = 171.0
$ws.Range("B4").Value = 17.0
$ws.Range("B5").Value = 67.0
$ws.Range("B6").Value = 48.0
$ws.Range("B7").Value = 10.0
$ws.Range("B8").Value = 29.0
